$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows before the current last (blank) row 17 so that the
# pre-existing blank row becomes row 23 (keeping its original formatting),
# and rows 17-22 are fresh rows for the new "employee" records.
$ws.Range("A17:A22").EntireRow.Insert()

# ---- Row 17: employee / S_EMPLOYEE ----
$ws.Range("A17").Value = "employee"
$ws.Range("B17").Value = "S_EMPLOYEE"
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = "PERNR"

# ---- Row 18: employee / S_PA0000 ----
$ws.Range("A18").Value = "employee"
$ws.Range("B18").Value = "S_PA0000"
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = "PERNR,ENDDA,BEGDA,MASSN,MASSG,STAT2"

# ---- Row 19: employee / S_PA0001 ----
$ws.Range("A19").Value = "employee"
$ws.Range("B19").Value = "S_PA0001"
$ws.Range("C19").Value = 4
$ws.Range("D19").Value = 3
$ws.Range("E19").Value = "PERNR,ENDDA,BEGDA,BUKRS,WERKS,VDSK1,BTRTL,KOSTL,KOKRS,PERSG,PERSK,ORGEH,OTYPE,MSTBR"

# ---- Row 20: employee / S_PA0002 ----
$ws.Range("A20").Value = "employee"
$ws.Range("B20").Value = "S_PA0002"
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = "PERNR,ENDDA,BEGDA,INITS,NACHN,NACH2,VORNA,TITEL,MIDNM,SPRSL"

# ---- Row 21: employee / S_PA0006 ----
$ws.Range("A21").Value = "employee"
$ws.Range("B21").Value = "S_PA0006"
$ws.Range("C21").Value = 4
$ws.Range("D21").Value = 4
$ws.Range("E21").Value = "SUBTY,ENDDA,BEGDA,ANSSA,STRAS,ORT01,ORT02,PSTLZ,LAND1,LOCAT,ADR03,ADR04,STATE,HSNMR,BLDNG,FLOOR,STRDS,COUNC,RCTVC,COM01,NUM01,COM02,NUM02,COM03,NUM03,COM04,NUM04,COM05,NUM05,COM06,NUM06"

# ---- Row 22: employee / S_PA0105 ----
$ws.Range("A22").Value = "employee"
$ws.Range("B22").Value = "S_PA0105"
$ws.Range("C22").Value = 4
$ws.Range("D22").Value = 5
$ws.Range("E22").Value = "PERNR,SUBTY,ENDDA,BEGDA,USRTY"

# ---- Row 23 (previously the blank placeholder row 17): employee / S_INFOTYPE_TEXT ----
$ws.Range("A23").Value = "employee"
$ws.Range("B23").Value = "S_INFOTYPE_TEXT"
$ws.Range("C23").Value = 4
$ws.Range("D23").Value = 6
$ws.Range("E23").Value = "INFTY,SUBTY,ENDDA,BEGDA"

# Formatting for the new rows 17-22: right-aligned, no border.
# Columns C/D get a thousands-separator number format; column E stays general text.
$ws.Range("C17:D22").NumberFormat = "#,##0"
$ws.Range("C17:D22").HorizontalAlignment = -4152
$ws.Range("E17:E22").HorizontalAlignment = -4152

# Row heights
$ws.Rows.Item(1).RowHeight = 19.5
for ($r = 2; $r -le 16; $r++) {
    $ws.Rows.Item($r).RowHeight = 18.75
}
for ($r = 17; $r -le 22; $r++) {
    $ws.Rows.Item($r).RowHeight = 18
}
$ws.Rows.Item(23).RowHeight = 19.5
